$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "mode" stats row into a new mode ("testb1"), same pattern
# as the existing test1/test2/test3 rows (label only, matching the diff).
$ws.Range("A7").Value = "testb1"

# Move the selection/cursor to reflect where the new data was entered
[void]$ws.Range("B7").Select()

# Restore the zoom level recorded for the sheet view
$win = $excel.ActiveWindow
$win.Zoom = 150
